{"js": "// Update the \"Liczba linii kodu\" (lines-of-code) column in the file-listing\n// table for the six rows whose counts grew, matching the source-repo upload.\n// Table layout (0-indexed columns: 0 = file name, 1 = size, 2 = line count):\n//   Row 2  Kolo.php        49 -> 54\n//   Row 3  Kwadrat.php     46 -> 50\n//   Row 4  Pieciokat.php   50 -> 55\n//   Row 5  Prostokat.php   50 -> 54\n//   Row 7  Szesciokat.php  46 -> 52\n//   Row 8  Trojkat.php     58 -> 61\n// (Row 6, Router.php, keeps its value of 33 and is intentionally skipped.)\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = {\n  2: \"54\",\n  3: \"50\",\n  4: \"55\",\n  5: \"54\",\n  7: \"52\",\n  8: \"61\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const row = Number(rowIndex);\n  // Replace the cell's text in place (rather than the `.value =` shortcut)\n  // so the existing run/paragraph markup is preserved untouched.\n  table.getCell(row, 2).getRange().insertText(updates[rowIndex], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Liczba linii kodu\" (lines-of-code) column in the file-listing\n# table for the six rows whose counts grew, matching the source-repo upload.\n# Table layout (col 1 = file name, col 2 = size, col 3 = line count):\n#   Row 3  Kolo.php        49 -> 54\n#   Row 4  Kwadrat.php     46 -> 50\n#   Row 5  Pieciokat.php   50 -> 55\n#   Row 6  Prostokat.php   50 -> 54\n#   Row 8  Szesciokat.php  46 -> 52\n#   Row 9  Trojkat.php     58 -> 61\n# (Row 7, Router.php, keeps its value of 33 and is intentionally skipped.)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    3 = \"54\"\n    4 = \"50\"\n    5 = \"55\"\n    6 = \"54\"\n    8 = \"52\"\n    9 = \"61\"\n}\n\nforeach ($row in $updates.Keys) {\n    $t.Cell($row, 3).Range.Text = $updates[$row]\n}\n"}
